# Renaming terms in tabular-defined schemas and table metadata, and
# removing the "!_Table of contents" sheet (examples/biochemical_models/template.xlsx)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Delete the "!_Table of contents" sheet entirely.
# ---------------------------------------------------------------------------
$toc = $wb.Worksheets.Item("!_Table of contents")
[void]$toc.Delete()

# ---------------------------------------------------------------------------
# Helper values
# ---------------------------------------------------------------------------
$oldDate = "2019-09-23 10:00:35"
$newDate = "2019-10-10 02:11:39"

# ---------------------------------------------------------------------------
# 2) "!_Schema" sheet
# ---------------------------------------------------------------------------
$schema = $wb.Worksheets.Item("!_Schema")

# Row-1 table metadata cell (locked cell -> unlock, edit, relock)
$a1 = $schema.Cells.Item(1, 1)
$a1.Locked = $false
$a1.Value = "!!ObjTables Type='Schema' Description='Table/model and column/attribute definitions' Date='" + $newDate + "' ObjTablesVersion='0.0.8'"
$a1.Locked = $true

# Column D ("!Type") attribute-type renames
for ($r = 4; $r -le 19; $r++) {
    $cell = $schema.Cells.Item($r, 4)
    $v = $cell.Text
    if ($v -eq "SlugAttribute") { $cell.Value = "Slug" }
    elseif ($v -eq "StringAttribute") { $cell.Value = "String" }
    elseif ($v -eq "BooleanAttribute") { $cell.Value = "Boolean" }
    elseif ($v -eq "ManyToOneAttribute") { $cell.Value = "ManyToOne" }
}

# ---------------------------------------------------------------------------
# 3) "!Compound" sheet
# ---------------------------------------------------------------------------
$compound = $wb.Worksheets.Item("!Compound")

$a1 = $compound.Cells.Item(1, 1)
$a1.Locked = $false
$a1.Value = "!!ObjTables Type='Data' Id='Compound' Description='Compound' Name='Compound' Date='" + $newDate + "' ObjTablesVersion='0.0.8'"
$a1.Locked = $true

$cmt = $compound.Range("A2").Comment
[void]$cmt.Text('Select a value from "!Model:1" or blank.')

$dv = $compound.Range("A3:A12").Validation
$dv.Modify(3, 2, 1, "'!Model'!`$B`$1:`$XFD`$1")
$dv.ErrorTitle = "Model"
$dv.ErrorMessage = 'Value must be a value from "!Model:1" or blank.'
$dv.InputTitle = "Model"
$dv.InputMessage = 'Select a value from "!Model:1" or blank.'

# ---------------------------------------------------------------------------
# 4) "!Model" sheet
# ---------------------------------------------------------------------------
$model = $wb.Worksheets.Item("!Model")

$a1 = $model.Cells.Item(1, 1)
$a1.Locked = $false
$a1.Value = "!!ObjTables Type='Data' Id='Model' Description='Model' Name='Model' Date='" + $newDate + "' ObjTablesVersion='0.0.8'"
$a1.Locked = $true

# ---------------------------------------------------------------------------
# 5) "!Reaction" sheet
# ---------------------------------------------------------------------------
$reaction = $wb.Worksheets.Item("!Reaction")

$a1 = $reaction.Cells.Item(1, 1)
$a1.Locked = $false
$a1.Value = "!!ObjTables Type='Data' Id='Reaction' Description='Reaction' Name='Reaction' Date='" + $newDate + "' ObjTablesVersion='0.0.8'"
$a1.Locked = $true

$cmt = $reaction.Range("A2").Comment
[void]$cmt.Text('Select a value from "!Model:1" or blank.')

$dv = $reaction.Range("A3:A12").Validation
$dv.Modify(3, 2, 1, "'!Model'!`$B`$1:`$XFD`$1")
$dv.ErrorTitle = "Model"
$dv.ErrorMessage = 'Value must be a value from "!Model:1" or blank.'
$dv.InputTitle = "Model"
$dv.InputMessage = 'Select a value from "!Model:1" or blank.'

Write-Host "edit complete"
